# Atualização das bases do grupo 13
# Adds the new quarter (01/10/2024) for each of the three regions
# (Brasil, Nordeste, Sergipe) present in the sheet, inserting each new
# row immediately after the existing rows of its region so the table
# stays grouped by region, then ordered by quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variable = "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência"

# --- Insert Brasil's new quarter right after the existing Brasil block ---
# Existing Brasil block occupies rows 2-24; inserting a row above row 25
# (the current first Nordeste row) pushes Nordeste/Sergipe down by one.
$ws.Rows.Item(25).Insert()
$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "Brasil"
$ws.Cells.Item(25, 2).Value = $variable
$ws.Cells.Item(25, 3).Value = "01/10/2024"
$ws.Cells.Item(25, 4).Value = 50.87

# --- Insert Nordeste's new quarter right after the existing Nordeste block ---
# After the previous insert, the Nordeste block now occupies rows 26-48,
# and Sergipe now starts at row 49. Inserting above row 49 pushes Sergipe
# down by one more row.
$ws.Rows.Item(49).Insert()
$ws.Cells.Item(49, 3).NumberFormat = "@"
$ws.Cells.Item(49, 1).Value = "Nordeste"
$ws.Cells.Item(49, 2).Value = $variable
$ws.Cells.Item(49, 3).Value = "01/10/2024"
$ws.Cells.Item(49, 4).Value = 44.25

# --- Append Sergipe's new quarter at the end of the table ---
# After the two inserts above, Sergipe occupies rows 50-72, so the new
# quarter goes in row 73.
$ws.Rows.Item(73).Insert()
$ws.Cells.Item(73, 3).NumberFormat = "@"
$ws.Cells.Item(73, 1).Value = "Sergipe"
$ws.Cells.Item(73, 2).Value = $variable
$ws.Cells.Item(73, 3).Value = "01/10/2024"
$ws.Cells.Item(73, 4).Value = 46.4

$wb.Save()
